$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("D18").Value = "['ELM-2NA-Instalções Elétricas', -]"
$ws.Range("E18").Value = "ELM-2NA-Máquinas Elétricas"
$ws.Range("F18").Value = "[-, 'ELM-2NA-Lab. De Máquinas elétricas']"

# Row 20
$ws.Range("C20").Value = "ELM-1NA-Circuitos Elétricos 1"
$ws.Range("D20").Value = "[-, 'ELM-2NA-Instalções Elétricas']"
$ws.Range("F20").Value = "[-, 'ELM-2NA-Instalções Elétricas']"

# Row 21
$ws.Range("C21").Value = "ELM-1NA-Circuitos Elétricos 1"
$ws.Range("D21").Value = "['ELM-2NA-Lab. De Máquinas elétricas', -]"
$ws.Range("E21").Value = "Allan Cupertino-Máquinas Elétricas"
$ws.Range("F21").Value = "[-, 'ELM-2NA-Instalções Elétricas']"
